# "remove column from alcohol data"
#
# The measurement sheet (Sheet1) has an extra column (M) whose values are a
# duplicate/obsolete copy of the data that should actually live one column
# over. The fix is to delete column M outright: the old column N slides
# left and becomes the new column M, and the sheet's used range shrinks
# from A1:N119 to A1:M119.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)

# Remove column M (13th column) from the data sheet; everything to the
# right (column N) shifts one column to the left.
$ws1.Range("M:M").EntireColumn.Delete()

# Reflect the edit the way a user would: zoom out a bit on every sheet and
# land the selection on the now-last data column.
$ws1.Activate()
$excel.ActiveWindow.Zoom = 85
$ws1.Range("M1").Select()

$ws2.Activate()
$excel.ActiveWindow.Zoom = 85

$ws3.Activate()
$excel.ActiveWindow.Zoom = 85

$ws1.Activate()
